# "Test für empty text" — fills in the previously-empty SCALA help text on
# Tabelle2, fixes the example filename/scale text on Tabelle1, and leaves
# the workbook with Tabelle2 as the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Tabelle2")

# --- Tabelle1: example AUDIO/SCALA values -------------------------------
# D4 held an example audio filename; switch the extension from wav to png
# (matches the rest of the template, which otherwise only ever shows .png).
$ws1.Range("D4").Value = "audio_file_name.png"

# --- Tabelle2: fill in the previously empty SCALA help text (B9) --------
$cell = $ws2.Range("B9")

$plainPart = "Benutze das, um eine Scala als Antwortmöglichkeit einzufügen. "
$redPart = "Achtung: der Text muss immer wie folgt aussehen: minimumZahl(text für minimum)-maximumZahl(text für maximum). Beispiel: 0(wenig motiviert)-5(sehr motiviert)"
$cell.Value = $plainPart + $redPart

$startPos = $plainPart.Length + 1
$redLen = $redPart.Length
$redRun = $cell.Characters($startPos, $redLen)
$redRun.Font.Size = 12
$redRun.Font.Color = 192
$redRun.Font.Name = "Calibri (Textkörper)"

# Row 9 now wraps across several lines like the other help rows in this
# column, so grow it to match.
$ws2.Rows.Item(9).RowHeight = 85

# J4 held the SCALA example text; the minimum label should start at 0, not 1.
$ws1.Range("J4").Value = "0(wenig motiviert)-5(sehr motiviert)"

# --- Window / selection state --------------------------------------------
# Move the active tab + selection from Tabelle1 to Tabelle2, leaving behind
# a selection on Tabelle1 at H6.
$ws1.Range("H6").Select()
$ws2.Select()
$ws2.Range("B7").Select()
